# chore: adapt column header formatting to respective input file names
# - rename the "_old"/"_new" header suffixes to "_FV2304"/"_FV2310"
# - turn the data range into a native Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRange = $ws.Range("A1:U1")

# Stash the header row's existing formatting (bold/fill/border/centered+wrap)
# in an unused scratch row so we can restore it after the table is created.
# (ListObjects.Add bakes in a header dxf when the header cells already carry
# custom formatting at creation time, so we temporarily clear it.)
$scratch = $ws.Range("A100:U100")
$headerRange.Copy()
$scratch.PasteSpecial(-4122)
$headerRange.ClearFormats()

# Rename the header cells to the new "<formatversion>" suffix scheme.
$ws.Range("A1").Value = "Segmentname_FV2304"
$ws.Range("B1").Value = "Segmentgruppe_FV2304"
$ws.Range("C1").Value = "Segment_FV2304"
$ws.Range("D1").Value = "Datenelement_FV2304"
$ws.Range("E1").Value = "Segment ID_FV2304"
$ws.Range("F1").Value = "Code_FV2304"
$ws.Range("G1").Value = "Qualifier_FV2304"
$ws.Range("H1").Value = "Beschreibung_FV2304"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("J1").Value = "Bedingung_FV2304"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2310"
$ws.Range("M1").Value = "Segmentgruppe_FV2310"
$ws.Range("N1").Value = "Segment_FV2310"
$ws.Range("O1").Value = "Datenelement_FV2310"
$ws.Range("P1").Value = "Segment ID_FV2310"
$ws.Range("Q1").Value = "Code_FV2310"
$ws.Range("R1").Value = "Qualifier_FV2310"
$ws.Range("S1").Value = "Beschreibung_FV2310"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("U1").Value = "Bedingung_FV2310"

# Turn A1:U53 into a native table.
$listObj = $ws.ListObjects.Add(1, $ws.Range("A1:U53"), 0, 1)
$listObj.Name = "Table1"

# Restore the header row's original formatting, then clean up the scratch row.
$scratch.Copy()
$headerRange.PasteSpecial(-4122)
$scratch.ClearContents()
$scratch.ClearFormats()

# Freeze the header row (split below row 1, bottom-left pane active).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
